# Flask PPTX-generator template: turn the static sample copy into a
# Jinja-style template by swapping the literal placeholder strings for
# {{token}} markers, and nudge the "titulo" textbox a bit further down.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)

# --- "titulo" textbox: move down (y: 895483 -> 925963 EMU) and templatize ---
$tituloShape = $s.Shapes.Item(6)
# Shape.Top is a COM Single (points); 72.91047668457031pt is the exact
# float32 value that round-trips to 925963 EMU (925963 / 12700).
$tituloShape.Top = 72.91047668457031
$tituloShape.TextFrame.TextRange.Text = "{{titulo}}"

# --- "data" textbox: just templatize ---
$s.Shapes.Item(7).TextFrame.TextRange.Text = "{{data}}"

# --- "link" textbox (the one under the first article): just templatize ---
$s.Shapes.Item(8).TextFrame.TextRange.Text = "{{link}}"

# --- "resumo" textbox: templatize, then restore the auto-fit height the ---
# --- engine recalculates after the run-length change (1060034 EMU). ---
$resumoShape = $s.Shapes.Item(18)
$resumoShape.TextFrame.TextRange.Paragraphs(1).TextRange.Text = "{{resumo}}"
$resumoShape.Height = 83.46724700927734
